$d = $word.ActiveDocument

# 1. Fix the typo "staffs" -> "staff's" (Word's AutoFormat turns the
#    straight apostrophe into a typographic one automatically).
$find1 = $d.Content
$find1.Find.Execute("staffs", $false, $false, $false, $false, $false, $true, 1, $false, "staff's", 2)

# 2. Split the run right after "staff's" so "n Backend Displays the
#    staff's" and " permission to login in the system." become two
#    separate runs (matching how the paragraph was re-authored).
$tail = $d.Content
$tail.Find.Execute(" permission to login in the system.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Bold = 1
$tail.Bold = 0

# 3. Relocate the "_GoBack" bookmark: it used to sit right after "The
#    Login i" (an artifact of the previous edit position); it should
#    now wrap the picture paragraph that immediately follows the text
#    paragraph we just edited.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$textRange = $d.Content
$textRange.Find.Execute("permission to login in the system", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$textPara = $textRange.Paragraphs(1)
$picPara = $textPara.Next()
$d.Bookmarks.Add("_GoBack", $picPara.Range)
